# #5: cash & deposit done
# Fill in the bank/deposit columns on the "存款" (deposit) sheet and append
# the standard metadata columns (property_category .. index) that the other
# asset sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# --- Header row (row 1) -----------------------------------------------
# B1/C1/D1 (bank/deposit_type/currency) already hold the right headers;
# E1 was "legislator_name"(邱志偉) by mistake - should be "owner", and F1
# held a stray numeric value instead of the "total" header.
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# New trailing metadata header cells G1:M1, matching the style of the
# existing header cells (bold, bordered, centered).
$ws.Range("B1").Copy()
$ws.Range("G1:M1").PasteSpecial(-4122)

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# --- Data rows (rows 2-9) ----------------------------------------------
# Correct bank name / deposit type per row (previously each row's B/C
# referenced the wrong shared string - effectively shifted by one entry).
$banks = @(
    "彰化商業銀行古亭分行",
    "臺灣銀行群賢分行",
    "高雄銀行市府分行",
    "華南商業銀行淡水分行",
    "中國信託商業銀行高雄分行",
    "中華郵政股份有限公司麟洛郵局",
    "中華郵政股份有限公司西甲郵局",
    "臺灣中小企業銀行前鎮分行"
)
$depositTypes = @(
    "活期儲蓄存款",
    "活期存款",
    "活期儲蓄存款",
    "活期儲蓄存款",
    "活期存款",
    "活期存款",
    "定期存款",
    "活期存款"
)

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $banks[$i]
    $ws.Cells.Item($r, 3).Value = $depositTypes[$i]
}

# New trailing metadata cells G2:M9, copying the data-row style (no
# border, regular weight) from an existing data cell.
$ws.Range("B2").Copy()
$ws.Range("G2:M9").PasteSpecial(-4122)

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2
    $idx = $i + 46
    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2012-03-06"
    $ws.Cells.Item($r, 10).Value = "邱志偉"
    $ws.Cells.Item($r, 11).Value = 1744
    $ws.Cells.Item($r, 12).Value = "tmpba221"
    $ws.Cells.Item($r, 13).Value = $idx
}

Write-Host "deposit sheet updated"
